$d = $word.ActiveDocument

function Insert-AsNewRun([string]$text) {
    # Range of paragraph 1, collapsed to just before the paragraph mark.
    $p = $d.Paragraphs(1).Range
    $insertStart = $p.End - 1
    $p.InsertAfter($text)

    # Wrapping the freshly-inserted span in a temporary bookmark (and then
    # removing the bookmark) forces the new text to stay in its own run
    # instead of being re-merged into the preceding run, even though the
    # two runs end up with identical (empty) formatting. That's what
    # produces the separate <w:r> elements seen in the target edit.
    $newRange = $d.Range($insertStart, $insertStart + $text.Length)
    $bmName = "TempSplitMark"
    $d.Bookmarks.Add($bmName, $newRange)
    $d.Bookmarks($bmName).Delete()
}

Insert-AsNewRun " ("
Insert-AsNewRun "Changed main"
Insert-AsNewRun ")"
